# Update IPC results for SP-D parallel execution
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value  = 2.184313988853944
$ws.Range("B3").Value  = 1.886159094700809
$ws.Range("B4").Value  = 2.161741777228566
$ws.Range("B5").Value  = 2.134971957482482
$ws.Range("B6").Value  = 2.265207983438145
$ws.Range("B7").Value  = 2.388832565585126
$ws.Range("B8").Value  = 1.832837670377314
$ws.Range("B9").Value  = 1.672161674023858
$ws.Range("B10").Value = 1.521127160740514
$ws.Range("B11").Value = 1.634850098626581
$ws.Range("B12").Value = 1.633760465320453
$ws.Range("B13").Value = 1.470056583249274
$ws.Range("B14").Value = 2.059414418759483
$ws.Range("B15").Value = 2.069838129682612
$ws.Range("B16").Value = 1.676471650687426
$ws.Range("B17").Value = 1.669251038708748
$ws.Range("B18").Value = 1.432430522853284
$ws.Range("B19").Value = 1.860164591717093
$ws.Range("B20").Value = 1.752263920085888
$ws.Range("B21").Value = 1.689741302980548
$ws.Range("B22").Value = 2.658605918231318
